$d = $word.ActiveDocument

# 1. Append " (60)" after "LED-Streifen mit ext. Netzteil"
$d.Content.Find.Execute("LED-Streifen mit ext. Netzteil", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LED-Streifen mit ext. Netzteil (60)", 2)

# 2. Update the cached date field results in the footers from 17.09.2021 to 24.09.2021
foreach ($sec in $d.Sections) {
    $footers = @($sec.Footers(1), $sec.Footers(2), $sec.Footers(3))
    foreach ($ftr in $footers) {
        if ($ftr -ne $null -and $ftr.Exists) {
            $ftr.Range.Find.Execute("17.09.2021", $false, $false, $false, $false, $false,
                                     $true, 1, $false, "24.09.2021", 2)
        }
    }
}
